$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper note: several Price values in column D look like plain decimal
# numbers (e.g. "516.81"); Excel's normal cell-value parser would silently
# convert those into numeric values. The source data keeps them as literal
# text (matching the other "thousand-dot" prices like "58.354.46"), so for
# those specific cells we force the Text number format before assigning the
# string value, cell by cell, so each one keeps its literal text content.

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "58.354.46"
$ws.Range("E2").Value = "  -1.45%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "2.571.91"
$ws.Range("E3").Value = "  -2.43%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.81"
$ws.Range("E5").Value = "  -2.26%  "

# Row 6 (Solana)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.61"
$ws.Range("E6").Value = "  -5.19%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.09%  "

# Row 8 (XRP)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  -1.32%  "

# Row 9 (LidoStakedEther)
$ws.Range("D9").Value = "2.585.37"
$ws.Range("E9").Value = "  -2.43%  "

# Row 10 (Toncoin)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  -2.75%  "

# Row 11 (Dogecoin)
$ws.Range("E11").Value = "  -4.31%  "

# Row 12 (Cardano)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.327"
$ws.Range("E12").Value = "  -2.62%  "

# Row 13 (TRON)
$ws.Range("E13").Value = "  +0.98%  "

# Row 14 (WrappedliquidstakedEther2.0)
$ws.Range("D14").Value = "3.028.07"
$ws.Range("E14").Value = "  -2.36%  "

# Row 15 (WrappedBTC)
$ws.Range("D15").Value = "58.337.38"
$ws.Range("E15").Value = "  -1.50%  "

# Row 16 (Avalanche)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.16"
$ws.Range("E16").Value = "  -2.53%  "

# Row 17 (WrappedEther)
$ws.Range("D17").Value = "2.572.43"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18 (ShibaInu)
$ws.Range("E18").Value = "  -3.87%  "

# Row 19 (BitcoinCash)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.29"

# Row 20 (Polkadot)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  -3.42%  "

# Row 21 (Chainlink)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.04"
$ws.Range("E21").Value = "  -4.66%  "

# Row 22 (Uniswap)
$ws.Range("E22").Value = "  +0.20%  "

# Row 23 (Dai)
$ws.Range("E23").Value = "  -0.11%  "

# Row 24 (Litecoin)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.73"
$ws.Range("E24").Value = "  -1.02%  "

# Row 25 (Kaspa)
$ws.Range("E25").Value = "  -1.23%  "

# Row 26 (Binance-PegBSC-USD)
$ws.Range("E26").Value = "  -0.26%  "

# Row 27 (Polygon)
$ws.Range("E27").Value = "  -3.40%  "

# Row 28 (InternetComputer(DFINITY))
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.94"

# Row 29 (USDe)
$ws.Range("E29").Value = "  -0.05%  "

# Row 30 (PEPE)
$ws.Range("D30").Value = "0.0₃0702"
$ws.Range("E30").Value = "  -11.36%  "

# Row 31 (Aptos)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.80"
$ws.Range("E31").Value = "  -7.69%  "

# Row 32 (EthereumClassic)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.59"
$ws.Range("E32").Value = "  -1.90%  "

# Row 33 (PancakeSwap)
$ws.Range("E33").Value = "  -3.55%  "

# Row 34 (Monero)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.40"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35 (NEARProtocol)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.88"
$ws.Range("E35").Value = "  -5.91%  "

# Row 36 (ImmutableX)
$ws.Range("E36").Value = "  -5.05%  "

# Row 37 (OKB)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.21"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38 (SuiNetwork)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.819"
$ws.Range("E38").Value = "  -3.84%  "

# Row 39 (Fetch.AI)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.819"
$ws.Range("E39").Value = "  -2.42%  "

# Row 40 (Stacks)
$ws.Range("E40").Value = "  -2.77%  "

# Row 41 (Filecoin)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  -4.04%  "

# Row 42 (FirstDigitalUSD)
$ws.Range("E42").Value = "  -0.11%  "

# Row 43 (Bittensor)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "269.09"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 (WhiteBITCoin)
$ws.Range("E44").Value = "  +0.14%  "

# Row 45 (Mantle)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.586"
$ws.Range("E45").Value = "  -1.96%  "

# Row 46 (Stellar)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0940"
$ws.Range("E46").Value = "  -3.77%  "

# Row 47 (Hedera)
$ws.Range("E47").Value = "  -3.52%  "

# Row 48 / 49 swap: EnergySwap and Maker rows exchange places (with updated
# values) in the ranking.
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.969.74"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.22"
$ws.Range("E49").Value = "  -4.94%  "

# Row 50 (VeChain)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0217"
$ws.Range("E50").Value = "  -4.81%  "

# Row 51 (RenderToken)
$ws.Range("E51").Value = "  -6.43%  "
